{"js": "// Replace the 100 arithmetic \"problem\" strings in the single 20x5 table\n// with their new values, in document (row-major) order. The mapping below\n// was derived from the canonical OOXML diff: it lists, in table order, the\n// expected \"before\" text for each cell followed by the \"after\" text that\n// must replace it. Matching on the expected \"before\" value (rather than\n// blindly overwriting every cell) makes the script a no-op / safe-guarded\n// if a cell has already been updated or the content doesn't line up.\nconst expected = [\n  \"29+66=95\", \"68+24=92\", \"8+57=65\", \"17+15=32\", \"5+39=44\",\n  \"42-24=18\", \"59+19=78\", \"6+28=34\", \"20-19=1\", \"83-38=45\",\n  \"70-64=6\", \"9+17=26\", \"45+17=62\", \"71-2=69\", \"51-34=17\",\n  \"34-29=5\", \"53-14=39\", \"43-15=28\", \"83-6=77\", \"56+25=81\",\n  \"62-36=26\", \"81-44=37\", \"44-6=38\", \"4+57=61\", \"9+56=65\",\n  \"97-89=8\", \"6+58=64\", \"17+57=74\", \"47+16=63\", \"58+26=84\",\n  \"17-9=8\", \"5+48=53\", \"66+8=74\", \"46-8=38\", \"70-56=14\",\n  \"69+3=72\", \"39+23=62\", \"81-56=25\", \"69+19=88\", \"61-52=9\",\n  \"73-15=58\", \"57+37=94\", \"46-18=28\", \"17+57=74\", \"53-14=39\",\n  \"17+78=95\", \"47+25=72\", \"29+43=72\", \"52-25=27\", \"92-4=88\",\n  \"72-17=55\", \"61-14=47\", \"61-13=48\", \"47+15=62\", \"69+7=76\",\n  \"81-77=4\", \"38+46=84\", \"64-25=39\", \"76-47=29\", \"96-9=87\",\n  \"63-6=57\", \"90-71=19\", \"17+27=44\", \"46+48=94\", \"75+6=81\",\n  \"63-49=14\", \"17+69=86\", \"73-25=48\", \"41-33=8\", \"13+48=61\",\n  \"28+47=75\", \"39+44=83\", \"33+48=81\", \"18+55=73\", \"72-8=64\",\n  \"18+45=63\", \"6+88=94\", \"36-9=27\", \"77+14=91\", \"59+27=86\",\n  \"95-56=39\", \"30-12=18\", \"26+9=35\", \"53-44=9\", \"19+64=83\",\n  \"43-29=14\", \"37-18=19\", \"8+75=83\", \"70-22=48\", \"91-72=19\",\n  \"38+6=44\", \"40-4=36\", \"7+67=74\", \"17+5=22\", \"19+19=38\",\n  \"35+8=43\", \"93-26=67\", \"69+19=88\", \"83-45=38\", \"9+53=62\"\n];\nconst replacement = [\n  \"70-64=6\", \"92-37=55\", \"33-25=8\", \"16+56=72\", \"92-37=55\",\n  \"36+45=81\", \"22+69=91\", \"80-45=35\", \"60-4=56\", \"24-19=5\",\n  \"55+38=93\", \"22-16=6\", \"90-53=37\", \"12+19=31\", \"74+19=93\",\n  \"71-63=8\", \"36+47=83\", \"18+15=33\", \"55-17=38\", \"19+24=43\",\n  \"43+48=91\", \"51-28=23\", \"24+67=91\", \"74-47=27\", \"32+49=81\",\n  \"82-3=79\", \"93-4=89\", \"19+13=32\", \"85-48=37\", \"80-45=35\",\n  \"38+3=41\", \"59+37=96\", \"37+29=66\", \"91-57=34\", \"61-27=34\",\n  \"16+55=71\", \"35+46=81\", \"26+36=62\", \"14+57=71\", \"23+69=92\",\n  \"85-36=49\", \"3+19=22\", \"84-78=6\", \"26+27=53\", \"25-16=9\",\n  \"45+7=52\", \"72-46=26\", \"65+16=81\", \"54+27=81\", \"40-11=29\",\n  \"11-4=7\", \"32-14=18\", \"61-49=12\", \"88-79=9\", \"29+47=76\",\n  \"2+59=61\", \"40-12=28\", \"24-17=7\", \"90-88=2\", \"67-19=48\",\n  \"85-37=48\", \"60-53=7\", \"31-6=25\", \"29+58=87\", \"65-6=59\",\n  \"92-25=67\", \"7+48=55\", \"47+4=51\", \"39+54=93\", \"27+5=32\",\n  \"45-39=6\", \"43+29=72\", \"82-48=34\", \"13+58=71\", \"57+25=82\",\n  \"95-86=9\", \"61-2=59\", \"37+19=56\", \"88-49=39\", \"24+7=31\",\n  \"24+67=91\", \"64-29=35\", \"62-38=24\", \"33-8=25\", \"50-27=23\",\n  \"54-35=19\", \"68+8=76\", \"62-6=56\", \"29+27=56\", \"74-39=35\",\n  \"45+29=74\", \"77+19=96\", \"45-29=16\", \"64-38=26\", \"48+24=72\",\n  \"80-67=13\", \"76-27=49\", \"85-26=59\", \"26+26=52\", \"26+27=53\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst current = table.values;\nconst newValues = [];\nlet cursor = 0;\nfor (let r = 0; r < current.length; r++) {\n  const row = current[r];\n  const newRow = [];\n  for (let c = 0; c < row.length; c++) {\n    const cellText = row[c];\n    if (cursor < expected.length && cellText === expected[cursor]) {\n      newRow.push(replacement[cursor]);\n    } else {\n      // Leave untouched if it doesn't match the expected \"before\" text\n      // (keeps the script safe against re-application / drift).\n      newRow.push(cellText);\n    }\n    cursor++;\n  }\n  newValues.push(newRow);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic \"problem\" strings in the single 20x5 table\n# with their new values, in document (row-major) order. The mapping below\n# was derived from the canonical OOXML diff: it lists, in table order, the\n# expected \"before\" text for each cell followed by the \"after\" text that\n# must replace it. Matching on the expected \"before\" value (rather than\n# blindly overwriting every cell) makes the script a no-op / safe-guarded\n# if a cell has already been updated or the content doesn't line up.\n$expected = @(\n  \"29+66=95\", \"68+24=92\", \"8+57=65\", \"17+15=32\", \"5+39=44\",\n  \"42-24=18\", \"59+19=78\", \"6+28=34\", \"20-19=1\", \"83-38=45\",\n  \"70-64=6\", \"9+17=26\", \"45+17=62\", \"71-2=69\", \"51-34=17\",\n  \"34-29=5\", \"53-14=39\", \"43-15=28\", \"83-6=77\", \"56+25=81\",\n  \"62-36=26\", \"81-44=37\", \"44-6=38\", \"4+57=61\", \"9+56=65\",\n  \"97-89=8\", \"6+58=64\", \"17+57=74\", \"47+16=63\", \"58+26=84\",\n  \"17-9=8\", \"5+48=53\", \"66+8=74\", \"46-8=38\", \"70-56=14\",\n  \"69+3=72\", \"39+23=62\", \"81-56=25\", \"69+19=88\", \"61-52=9\",\n  \"73-15=58\", \"57+37=94\", \"46-18=28\", \"17+57=74\", \"53-14=39\",\n  \"17+78=95\", \"47+25=72\", \"29+43=72\", \"52-25=27\", \"92-4=88\",\n  \"72-17=55\", \"61-14=47\", \"61-13=48\", \"47+15=62\", \"69+7=76\",\n  \"81-77=4\", \"38+46=84\", \"64-25=39\", \"76-47=29\", \"96-9=87\",\n  \"63-6=57\", \"90-71=19\", \"17+27=44\", \"46+48=94\", \"75+6=81\",\n  \"63-49=14\", \"17+69=86\", \"73-25=48\", \"41-33=8\", \"13+48=61\",\n  \"28+47=75\", \"39+44=83\", \"33+48=81\", \"18+55=73\", \"72-8=64\",\n  \"18+45=63\", \"6+88=94\", \"36-9=27\", \"77+14=91\", \"59+27=86\",\n  \"95-56=39\", \"30-12=18\", \"26+9=35\", \"53-44=9\", \"19+64=83\",\n  \"43-29=14\", \"37-18=19\", \"8+75=83\", \"70-22=48\", \"91-72=19\",\n  \"38+6=44\", \"40-4=36\", \"7+67=74\", \"17+5=22\", \"19+19=38\",\n  \"35+8=43\", \"93-26=67\", \"69+19=88\", \"83-45=38\", \"9+53=62\"\n)\n$replacement = @(\n  \"70-64=6\", \"92-37=55\", \"33-25=8\", \"16+56=72\", \"92-37=55\",\n  \"36+45=81\", \"22+69=91\", \"80-45=35\", \"60-4=56\", \"24-19=5\",\n  \"55+38=93\", \"22-16=6\", \"90-53=37\", \"12+19=31\", \"74+19=93\",\n  \"71-63=8\", \"36+47=83\", \"18+15=33\", \"55-17=38\", \"19+24=43\",\n  \"43+48=91\", \"51-28=23\", \"24+67=91\", \"74-47=27\", \"32+49=81\",\n  \"82-3=79\", \"93-4=89\", \"19+13=32\", \"85-48=37\", \"80-45=35\",\n  \"38+3=41\", \"59+37=96\", \"37+29=66\", \"91-57=34\", \"61-27=34\",\n  \"16+55=71\", \"35+46=81\", \"26+36=62\", \"14+57=71\", \"23+69=92\",\n  \"85-36=49\", \"3+19=22\", \"84-78=6\", \"26+27=53\", \"25-16=9\",\n  \"45+7=52\", \"72-46=26\", \"65+16=81\", \"54+27=81\", \"40-11=29\",\n  \"11-4=7\", \"32-14=18\", \"61-49=12\", \"88-79=9\", \"29+47=76\",\n  \"2+59=61\", \"40-12=28\", \"24-17=7\", \"90-88=2\", \"67-19=48\",\n  \"85-37=48\", \"60-53=7\", \"31-6=25\", \"29+58=87\", \"65-6=59\",\n  \"92-25=67\", \"7+48=55\", \"47+4=51\", \"39+54=93\", \"27+5=32\",\n  \"45-39=6\", \"43+29=72\", \"82-48=34\", \"13+58=71\", \"57+25=82\",\n  \"95-86=9\", \"61-2=59\", \"37+19=56\", \"88-49=39\", \"24+7=31\",\n  \"24+67=91\", \"64-29=35\", \"62-38=24\", \"33-8=25\", \"50-27=23\",\n  \"54-35=19\", \"68+8=76\", \"62-6=56\", \"29+27=56\", \"74-39=35\",\n  \"45+29=74\", \"77+19=96\", \"45-29=16\", \"64-38=26\", \"48+24=72\",\n  \"80-67=13\", \"76-27=49\", \"85-26=59\", \"26+26=52\", \"26+27=53\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellText = $cell.Range.Text\n    # Cell.Range.Text includes the trailing cell-mark (CR + BEL); strip it\n    # before comparing against the expected \"before\" value.\n    $cellText = $cellText.TrimEnd([char]13, [char]7)\n    if ($idx -lt $expected.Length -and $cellText -eq $expected[$idx]) {\n      $cell.Range.Text = $replacement[$idx]\n    }\n    $idx++\n  }\n}\n"}
